# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT even when it looks like a pure
# number (e.g. "601.15"), then restore the cell to the default "Normal" style
# so no stray number-format is left applied to the cell.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "66.766.06"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.499.72"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "601.15"
$ws.Range("E5").Value = "  -0.63%  "
Set-TextValue "D6" "147.42"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("D7").Value = "3.497.45"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.09%  "
Set-TextValue "D10" "0.142"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +5.97%  "
Set-TextValue "D12" "0.422"
$ws.Range("E12").Value = "  -1.31%  "
Set-TextValue "D13" "0.0000213"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "4.091.90"
$ws.Range("E14").Value = "  -0.19%  "
Set-TextValue "D15" "31.19"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("D16").Value = "3.500.49"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "66.782.10"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +7.56%  "
$ws.Range("E20").Value = "  -1.60%  "
Set-TextValue "D21" "15.35"
$ws.Range("E21").Value = "  -0.32%  "
Set-TextValue "D22" "433.74"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  -2.65%  "
Set-TextValue "D24" "79.57"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "3.638.61"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("E26").Value = "  +0.03%  "
Set-TextValue "D27" "0.0000119"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E29").Value = "  -5.89%  "
Set-TextValue "D30" "2.49"
$ws.Range("E30").Value = "  +0.12%  "
Set-TextValue "D31" "1.61"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  +0.32%  "
Set-TextValue "D33" "0.166"
$ws.Range("E33").Value = "  -0.96%  "
Set-TextValue "D34" "25.40"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "3.493.72"
$ws.Range("E35").Value = "  -0.17%  "
Set-TextValue "D36" "5.92"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("E37").Value = "  -2.76%  "
Set-TextValue "D38" "7.99"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  -0.01%  "
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +0.76%  "
Set-TextValue "D42" "169.99"
$ws.Range("E42").Value = "  -2.89%  "
Set-TextValue "D43" "2.09"
$ws.Range("E43").Value = "  -8.77%  "
$ws.Range("E44").Value = "  -0.20%  "
Set-TextValue "D45" "0.896"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D46" "1.34"
$ws.Range("E46").Value = "  +4.03%  "
Set-TextValue "D47" "45.83"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "28.57"
$ws.Range("E48").Value = "  -4.00%  "
Set-TextValue "D49" "7.47"
$ws.Range("E49").Value = "  -1.69%  "
Set-TextValue "D50" "2.42"
$ws.Range("E50").Value = "  -3.15%  "
Set-TextValue "D51" "0.970"
$ws.Range("E51").Value = "  -0.95%  "
